$wb = $excel.ActiveWorkbook

# --- Rotate the three "Acc_Upfront*" sheet names ---------------------------
# Before: Acc_Upfront3, Acc_Upfront1, Acc_Upfront2  (in tab order)
# After:  Acc_Upfront1, Acc_Upfront2, Acc_Upfront3
# Use temporary names to avoid name collisions mid-rotation, and rename back
# to the final names while keeping sheets in their existing tab positions.
$wb.Worksheets.Item("Acc_Upfront3").Name = "__Tmp_Upfront_A__"
$wb.Worksheets.Item("Acc_Upfront1").Name = "__Tmp_Upfront_B__"
$wb.Worksheets.Item("Acc_Upfront2").Name = "__Tmp_Upfront_C__"

$wb.Worksheets.Item("__Tmp_Upfront_A__").Name = "Acc_Upfront1"
$wb.Worksheets.Item("__Tmp_Upfront_B__").Name = "Acc_Upfront2"
$wb.Worksheets.Item("__Tmp_Upfront_C__").Name = "Acc_Upfront3"

# --- Update remembered selections on a couple of sheets ---------------------
$wsTransactions = $wb.Worksheets.Item("Transactions")
$wsTransactions.Activate() | Out-Null
$wsTransactions.Range("D5").Select() | Out-Null

# This is the sheet that used to be "Acc_Upfront2" (last tab, was the
# active/selected tab in the workbook) and is now named "Acc_Upfront3".
$wsUpfront3 = $wb.Worksheets.Item("Acc_Upfront3")
$wsUpfront3.Activate() | Out-Null
$wsUpfront3.Range("F18").Select() | Out-Null
